$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (bs_auto1)
$ws.Range("B2").Value = 41834364.20349601
$ws.Range("C2").Value = 1627522.883479001
$ws.Range("D2").Value = 924604.9999999969
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 6271

# Update row 3 (bs_auto2)
$ws.Range("B3").Value = 640861.9476030008
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 213620.6492010003
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 964

# Update row 4 (bs_auto3)
$ws.Range("B4").Value = 657615101.5700005
$ws.Range("C4").Value = 750500
$ws.Range("D4").Value = 106515
$ws.Range("E4").Value = 127856.7113140001
$ws.Range("F4").Value = 1501

# Delete rows 5-8
$ws.Range("A5:F8").EntireRow.Delete()
